# Fruta / hortaliza, semanal
# Insert a new weekly data row for "Feria Lagunitas de Puerto Montt - Durazno"
# right above the current row 84, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 84 (existing rows 84..186 shift down to 85..187)
$ws.Rows(84).Insert()

# Populate the constant / carried-over columns for the new row (same as the
# rest of the dataset for this market/product).
$ws.Range("A84").Value = 4
$ws.Range("B84").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C84").Value = "Los Lagos"
$ws.Range("D84").Value = 44579
$ws.Range("E84").Value = 10
$ws.Range("F84").Value = "Fruta"
$ws.Range("G84").Value = 100103
$ws.Range("H84").Value = "Frutos de hueso (carozo)"
$ws.Range("I84").Value = 100103004
$ws.Range("J84").Value = "Durazno"

# New row's specific data values
$ws.Range("K84").Value = "Carson"
$ws.Range("L84").Value = "Primera"
$ws.Range("M84").Value = 600
$ws.Range("N84").Value = 17000
$ws.Range("O84").Value = 18000
$ws.Range("P84").Value = 17500
$ws.Range("Q84").Value = "$/caja 15 kilos empedrada"
$ws.Range("R84").Value = "Región de O'Higgins"
$ws.Range("S84").Value = 1167
$ws.Range("T84").Value = 15
